# Apply a circular rotation of the data rows 2-4:
# new row2 = old row3, new row3 = old row4, new row4 = old row2
# Only columns D, J, K, L, M, O, P change; the rest stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "old" values for the columns that change, before overwriting anything.
$oldD = @{}
$oldJ = @{}
$oldK = @{}
$oldL = @{}
$oldM = @{}
$oldO = @{}
$oldP = @{}

foreach ($r in 2..4) {
    $oldD[$r] = $ws.Range("D$r").Value2
    $oldJ[$r] = $ws.Range("J$r").Value2
    $oldK[$r] = $ws.Range("K$r").Value2
    $oldL[$r] = $ws.Range("L$r").Value2
    $oldM[$r] = $ws.Range("M$r").Value2
    $oldO[$r] = $ws.Range("O$r").Value2
    $oldP[$r] = $ws.Range("P$r").Value2
}

# Mapping: target row -> source row (circular shift up: 2<-3, 3<-4, 4<-2)
$srcFor = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($r in 2..4) {
    $src = $srcFor[$r]
    $ws.Range("D$r").Value2 = $oldD[$src]
    $ws.Range("J$r").Value2 = $oldJ[$src]
    $ws.Range("K$r").Value2 = $oldK[$src]
    $ws.Range("L$r").Value2 = $oldL[$src]
    $ws.Range("M$r").Value2 = $oldM[$src]
    $ws.Range("O$r").Value2 = $oldO[$src]
    $ws.Range("P$r").Value2 = $oldP[$src]
}
